# Updating WeightMeasurements file 30/04/2018
# Adds a new measurement row (row 19) below the existing data (row 18),
# extends the two "fill-down" formula columns (F: Gain/Loss, G: BMI) to
# cover the new row, and updates the active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 19, copying the formatting down from row 18
# (this keeps every column's number format/style identical to the row
# above it, exactly like dragging the fill handle / pressing Enter at
# the bottom of the table would in the real Excel UI).
$ws.Rows("19:19").Insert(-4121, 0) | Out-Null   # xlShiftDown, xlFormatFromLeftOrAbove

# New data values (WeightID 18, dated 30/04/2018 = serial 43220)
$ws.Range("A19").Value2 = 18
$ws.Range("B19").Value2 = 43220
$ws.Range("C19").Value2 = 14.13
$ws.Range("D19").Value2 = 94.8
$ws.Range("E19").Value2 = 209
$ws.Range("H19").Value2 = 20.6

# Re-apply the two calculated columns across the whole column (existing
# rows keep their same formula/result, new row 19 now gets one too) so
# the fill-down formulas cover the new row.
$ws.Range("F3:F19").Formula = "=E3-E2"
$ws.Range("G2:G19").Formula = "=ROUND((D2/1.88)/1.88,2)"

# Update the view: scroll so row 4 is at the top and select E22, matching
# where the user left the sheet after adding the new entry.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E22").Select() | Out-Null
